$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.536.98"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "1.912.00"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.68%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.47"
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4827"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4076"
$ws.Range("E8").Value = "  -0.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08161"
$ws.Range("E9").Value = "  +1.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.011"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.45"
$ws.Range("E11").Value = "  +4.75%  "

$ws.Range("D12").Value = "1.933.48"
$ws.Range("E12").Value = "  +2.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.020"
$ws.Range("E13").Value = "  +1.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.111"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.41"
$ws.Range("E15").Value = "  +0.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06785"
$ws.Range("E16").Value = "  +2.68%  "

$ws.Range("E17").Value = "  +0.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001041"
$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.71"
$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("E20").Value = "  +0.51%  "

$ws.Range("D21").Value = "29.557.72"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.618"
$ws.Range("E22").Value = "  +1.18%  "

$ws.Range("E23").Value = "  +2.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.170"
$ws.Range("E24").Value = "  -1.50%  "

$ws.Range("D25").Value = "2.149.88"
$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.67"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.09"
$ws.Range("E27").Value = "  +1.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.324"
$ws.Range("E28").Value = "  +9.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.103"
$ws.Range("E29").Value = "  -1.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.74"
$ws.Range("E30").Value = "  +1.60%  "

$ws.Range("E31").Value = "  -2.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09567"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.542"
$ws.Range("E33").Value = "  +2.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.560"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.394"
$ws.Range("E35").Value = "  -2.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02268"
$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06113"
$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.176"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("E39").Value = "  +6.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5945"
$ws.Range("E40").Value = "  +1.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.936"
$ws.Range("E41").Value = "  -4.87%  "

$ws.Range("E42").Value = "  +0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.458"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.281"
$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07733"
$ws.Range("E45").Value = "  -3.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.45"
$ws.Range("E46").Value = "  +2.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5572"
$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.952"
$ws.Range("E48").Value = "  +1.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "115.22"
$ws.Range("E49").Value = "  +1.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.76"
$ws.Range("E50").Value = "  +1.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.053"
$ws.Range("E51").Value = "  +1.88%  "
